# edit.ps1 - apply "Few fixes on family, still need to work on decendants" commit
#
# Summary of required changes (see task diff):
#  Para 6  (5 points / commits):     yellow -> darkGreen; merge trailing "." run into main text run
#  Para 7  (5 points / validation):  add highlight yellow (was unhighlighted);
#                                    move the "_GoBack" bookmark here (to paragraph start)
#  Para 8  (10 points / send copy):  yellow -> darkGreen (leave the trailing, unformatted
#                                    space run at the end of the paragraph untouched)
#  Para 11 (15 points / look up info): yellow -> darkGreen across all of its runs
#  Para 13 (20 points / family):     yellow -> darkGreen; the "_GoBack" bookmark that used
#                                    to sit here moves to paragraph 7 (handled above)

$d = $word.ActiveDocument

# wdColorIndex constants used below
$wdYellow = 7
$wdGreen  = 11   # maps to the OOXML "darkGreen" highlight value

function Set-ParagraphHighlight($paraIndex, $colorIndex) {
    $p = $d.Paragraphs($paraIndex)
    # Exclude the trailing paragraph mark so we never spill over into the next paragraph.
    $full = $d.Range($p.Range.Start, $p.Range.End - 1)
    # Use Font.HighlightColorIndex (not Range.HighlightColorIndex) so that only the
    # characters actually covered by this range are recolored, not the whole paragraph.
    $full.Font.HighlightColorIndex = $colorIndex
}

# --- Paragraph 6: "(5 points): ... consistent commits ..." -----------------
# First, merge the unformatted trailing "." run into the preceding sentence run so the
# paragraph ends up with exactly two runs, matching the target structure. A no-op
# rewrite doesn't force a merge, so nudge the text by one extra character and then
# restore it; that does force the engine to coalesce the run.
$rngCommit = $d.Content
$rngCommit.Find.Execute("As a developer, I want to make consistent commits with good, descriptive messages")
$mergeRange = $d.Range($rngCommit.Start, $rngCommit.End + 1)
$origText = $mergeRange.Text
$mergeRange.Text = $origText + " "
$mergeRange2 = $d.Range($rngCommit.Start, $rngCommit.End + 2)
$mergeRange2.Text = $origText

Set-ParagraphHighlight 6 $wdGreen

# --- Paragraph 7: "(5 points): ... run validation ..." ---------------------
Set-ParagraphHighlight 7 $wdYellow

# Move the "_GoBack" bookmark to the start of this paragraph (adding a bookmark with an
# existing name relocates it, automatically removing it from its previous location).
$p7 = $d.Paragraphs(7)
$p7Start = $d.Range($p7.Range.Start, $p7.Range.Start)
$d.Bookmarks.Add("_GoBack", $p7Start)

# --- Paragraph 8: "(10 points): ... send a copy ..." ------------------------
Set-ParagraphHighlight 8 $wdGreen

# --- Paragraph 11: "(15 points): ... look up someone's information ..." ----
Set-ParagraphHighlight 11 $wdGreen

# --- Paragraph 13: "(20 points): ... immediate family ..." -----------------
Set-ParagraphHighlight 13 $wdGreen
